# Refresh the cryptocurrency snapshot values (prices / 1h volume moves)
# exported by the scheduled GitHub Actions job. Values are written as
# quote-prefixed text (leading "'") so Excel keeps them as plain strings
# -- matching the original inlineStr cell contents -- instead of coercing
# number-like text (e.g. "1.00", "0.998") into numeric values. The style
# is then reset to "Normal" so the quote-prefix formatting Excel applies
# behind the scenes doesn't linger on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = "'62.560.37"
    'E2' = "'  -3.56%  "
    'D3' = "'3.038.48"
    'E3' = "'  -3.02%  "
    'E4' = "'  -0.02%  "
    'D5' = "'539.60"
    'E5' = "'  -4.06%  "
    'E6' = "'  -10.30%  "
    'E7' = "'  +0.16%  "
    'D8' = "'3.027.46"
    'E8' = "'  -3.17%  "
    'D9' = "'0.484"
    'E9' = "'  -2.54%  "
    'D10' = "'6.35"
    'E10' = "'  -8.08%  "
    'E11' = "'  -2.86%  "
    'D12' = "'0.453"
    'E12' = "'  -1.78%  "
    'D13' = "'34.24"
    'E13' = "'  -4.15%  "
    'E14' = "'  -4.27%  "
    'D15' = "'3.530.81"
    'D16' = "'62.586.14"
    'E16' = "'  -3.54%  "
    'E17' = "'  -1.78%  "
    'D18' = "'3.034.59"
    'E18' = "'  -3.17%  "
    'D19' = "'6.52"
    'E19' = "'  -2.61%  "
    'D20' = "'476.91"
    'E20' = "'  -9.08%  "
    'D21' = "'13.21"
    'E21' = "'  -4.08%  "
    'D22' = "'0.688"
    'E22' = "'  -1.34%  "
    'D23' = "'6.96"
    'E23' = "'  -5.54%  "
    'D24' = "'76.72"
    'E24' = "'  -2.12%  "
    'D25' = "'12.01"
    'E25' = "'  -5.12%  "
    'D26' = "'0.998"
    'E26' = "'  +0.06%  "
    'E27' = "'  -3.86%  "
    'D28' = "'8.16"
    'E28' = "'  -4.88%  "
    'D29' = "'0.998"
    'E29' = "'  -0.30%  "
    'E30' = "'  -9.90%  "
    'D31' = "'25.88"
    'E31' = "'  -0.62%  "
    'E32' = "'  -3.52%  "
    'D33' = "'59.61"
    'E33' = "'  +12.98%  "
    'D34' = "'2.45"
    'E34' = "'  -6.75%  "
    'D35' = "'511.91"
    'E35' = "'  -7.98%  "
    'D36' = "'5.82"
    'E36' = "'  -2.75%  "
    'D37' = "'5.05"
    'E37' = "'  -5.72%  "
    'D38' = "'0.0394"
    'E38' = "'  -9.19%  "
    'D39' = "'3.043.42"
    'D40' = "'0.0777"
    'E40' = "'  -4.18%  "
    'D41' = "'0.116"
    'E41' = "'  -3.58%  "
    'D42' = "'7.95"
    'E42' = "'  -3.08%  "
    'E43' = "'  -9.21%  "
    'E44' = "'  +0.02%  "
    'D45' = "'0.248"
    'E45' = "'  -2.41%  "
    'E46' = "'  -6.82%  "
    'D47' = "'119.78"
    'E47' = "'  +0.62%  "
    'D48' = "'23.85"
    'E48' = "'  -4.11%  "
    'D49' = "'0.106"
    'E49' = "'  -2.40%  "
    'B50' = "'CoreDAO"
    'C50' = "'https://coinranking.com/coin/HFvoXUQh4+coredao-core"
    'D50' = "'2.32"
    'E50' = "'  +60.85%  "
    'B51' = "'PEPE"
    'C51' = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
    'D51' = "'0.0₃0485"
    'E51' = "'  -6.46%  "
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.Value = $updates[$cellRef]
    $cell.Style = 'Normal'
}
